$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Correct URL" column (B) for rows 9-29 that were previously blank
# (or, for a few rows, set equal to the Failed URL in column A, i.e. confirmed correct)
$ws.Range("B9").Value = 'http://www.mendeley.com/c/4974603702/g/2058663/smith-2006-the-vertebrate-fauna-of-ichauway-baker-county-ga/'
$ws.Range("B10").Value = 'not sure how to find this one'
$ws.Range("B11").Value = 'http://www.mendeley.com/c/5076417672/g/2058663/euliss-1996-ecological-studies-at-the-woodworth-study-area-terrestrial-bird-communities-on-the-woodworth-study-area-duplicated-copy-for-dcfs/'
$ws.Range("B12").Value = 'http://www.mendeley.com/c/5007365232/g/2058663/north-sterling-state-park-2012-north-sterling-state-park-birders-complete-checklist/'
$ws.Range("B13").Value = 'http://www.mendeley.com/c/5007859122/g/2058663/meyer-1985-classification-of-native-vegetation-at-the-woodworth-station-north-dakota/'
$ws.Range("B14").Value = 'http://www.mendeley.com/c/5007912492/g/2058663/higgins-1992-waterfowl-production-on-the-woodworth-station-in-south-central-north-dakota--1965-1981/'
$ws.Range("B15").Value = 'http://www.mendeley.com/c/4980636152/g/2058663/drew-1998-the-vascular-flora-of-ichauway--baker-county--georgia--a-remnant-longleaf-pine--wiregrass-ecosystem/'
$ws.Range("B16").Value = 'http://www.mendeley.com/c/5014915452/g/2058663/beckett-1982-forest-vegegation-and-vascular-flora-of-reed-brake-research-natural-area-alabama/'
$ws.Range("B17").Value = 'http://www.mendeley.com/c/5007859122/g/2058663/meyer-1985-classification-of-native-vegetation-at-the-woodworth-station-north-dakota/'
$ws.Range("B18").Value = 'http://www.mendeley.com/c/5007859122/g/2058663/meyer-1985-classification-of-native-vegetation-at-the-woodworth-station-north-dakota/'
$ws.Range("B19").Value = 'http://www.mendeley.com/c/5076426772/g/2058663/euliss-1996-ecological-studies-at-the-woodworth-study-area-upland-vegetation-at-the-woodworth-study-area/'
$ws.Range("B20").Value = 'http://www.mendeley.com/c/5076428272/g/2058663/euliss-1996-ecological-studies-at-the-woodworth-study-area-effects-of-water-level-changes-on-prairie-pothole-vegetation-structure-and-diversity-in-the-woodworth-study-area--north-dakota-duplicated-copy-for-dcfs/'
$ws.Range("B21").Value = 'http://www.mendeley.com/c/5018302512/g/2058663/shears-1999-central-arizona--phoenix-lter-deb-9714833-land-use-change-and-ecological-processes-in-an-urban-ecosystem-of-the-sonoran-desert-annual-progress-report-1999-2000/'
$ws.Range("B22").Value = 'http://www.mendeley.com/c/5009858662/g/2058663/hanson-1989-coleoptera-species-inhabiting-prairie-wetlands-of-the-cottonwood-lake-area-stutsman-county-north-dakota/'
$ws.Range("B23").Value = 'http://www.mendeley.com/c/5050303512/g/2058663/rice-2010-niche-preference-of-a-coprophagous-scarab-beetle--coleoptera--scarabaeidae--for-summer-moose-dung-in-denali-national-park--alaska/'
$ws.Range("B24").Value = 'http://www.mendeley.com/c/5017520222/g/2058663/cavey-2004-survey-report-on-the-leaf-beetles-of-cove-point-lng-property-and-vicinity-calvert-county-maryland/'
$ws.Range("B25").Value = 'http://vectormap.nhm.ku.edu/vectormap/'
$ws.Range("B26").Value = 'could not find'
$ws.Range("B27").Value = 'http://www.mendeley.com/c/5018230942/g/2058663/alabama-department-of-conservation-2000-outdoor-alabama-volumes-72-73/'
$ws.Range("B28").Value = 'http://www.mendeley.com/c/5001218602/g/2058663/genet-2001-the-lizard-community-of-a-subtropical-dry-forest--guanica-forest--puerto-rico/'
$ws.Range("B29").Value = 'could not find '

# Reflect the saved selection state (cursor left on the first empty row below the data)
$ws.Range("A35").Select() | Out-Null
